# Edit slide 74 ("Method parseVariableExpr() (continued)") of the
# "06 - Syntax Analysis" deck:
#   * Title: rename the method from parseVariableExpr to parseVariableCommon.
#   * Body : rewrite the while-loop body so that the leftBracket branch is
#            inlined (match/parseExpression/match instead of a single
#            parseIndexExpr() call) and the dot branch is inlined too
#            (match/match instead of a single parseFieldExpr() call).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(74)

# ---------------------------------------------------------------------
# Title placeholder: "Method parseVariableExpr()" -> "Method parseVariableCommon()"
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$titleTr = $title.TextFrame.TextRange
$titleText = $titleTr.Text
$oldName = "parseVariableExpr"
$newName = "parseVariableCommon"
$idx = $titleText.IndexOf($oldName)
if ($idx -ge 0) {
    $titleTr.Characters($idx + 1, $oldName.Length).Text = $newName
}

# ---------------------------------------------------------------------
# Content placeholder: rebuild the code listing paragraph by paragraph,
# each paragraph being a list of runs (so that the same text carries the
# same run-split boundaries as the final deck).
# ---------------------------------------------------------------------
$content = $s.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

$paraRuns = @(
  ,@("    while (", "scanner.getSymbol", "().", "isSelectorStarter", "())")
  ,@("      {")
  ,@("        if (", "scanner.getSymbol", "() == ", "Symbol.leftBracket", ")")
  ,@("          {")
  ,@("            match(", "Symbol.leftBracket", ");")
  ,@("            ", "parseExpression", "();")
  ,@("            match(", "Symbol.rightBracket", ");")
  ,@("          }")
  ,@("        else if (", "scanner.getSymbol", "() == Symbol.dot)")
  ,@("          {")
  ,@("            match(Symbol.dot);")
  ,@("            match(", "Symbol.identifier", ");")
  ,@("          }")
  ,@("      }")
  ,@("  }")
)

$cr = [char]13

# Build the plain paragraph texts (one string per paragraph) and push them
# all at once -- this lets the host line up each paragraph's pPr (indent /
# spacing / bullet suppression) positionally with the existing paragraphs.
$paraTexts = @()
foreach ($runList in $paraRuns) {
    $paraTexts += [string]::Join("", $runList)
}
$tr.Text = [string]::Join($cr, $paraTexts)

# Now re-carve the runs inside each paragraph so the run boundaries match
# the source runs (keeps independent formatting ranges available, e.g. for
# the identifiers vs. the surrounding punctuation/keywords).
for ($i = 0; $i -lt $paraRuns.Count; $i++) {
    $runList = $paraRuns[$i]
    if ($runList.Count -le 1) {
        continue
    }
    $paraStart = $tr.Paragraphs($i + 1, 1).Start
    $offset = 0
    foreach ($runText in $runList) {
        $len = $runText.Length
        if ($len -gt 0) {
            $tr.Characters($paraStart + $offset, $len).Text = $runText
        }
        $offset += $len
    }
}
